# Add a new data row (row 8) to Sheet1, mirroring the formatting of the
# preceding data row (row 7) and populating it with the 26-09-2025 gold
# price entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles/borders) from the last existing data row down to
# the new row so A8/B8 pick up the same cell styles as A7/B7.
$ws.Range("A7:B7").Copy()
$ws.Range("A8:B8").PasteSpecial(-4122)

# Now set the actual cell contents for the new row.
$ws.Range("A8").Value = "26-09-2025"
$ws.Range("B8").Value = "The price of gold in India today is ₹11,488 per gram for 24 karat gold, ₹10,530 per gram for 22 karat gold and ₹8,616 per gram for 18 karat gold (also called 999 gold)."
